# Sprint backlog update: add "Sprint 7 (M7)" tab with the M7 task list,
# mirroring the layout/format of the existing Sprint 6 (M6) sheet, and
# move the active-sheet/selection focus from Sprint 6 onto the new sheet.

$wb = $excel.ActiveWorkbook

# --- Deselect the previously active sheet (Sprint 6) and reset its selection ---
$ws6 = $wb.Worksheets.Item("Sprint 6 (M6)")
$ws6.Activate()
$ws6.Range("A1:G1").Select()

# --- Insert the new sheet right after "Sprint 6 (M6)" (i.e. at the end) ---
$lastIndex = $wb.Worksheets.Count
$ws6Again = $wb.Worksheets.Item($lastIndex)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws6Again)
$newSheet.Name = "Sprint 7 (M7)"

# --- Column widths matching the sprint-backlog layout ---
$newSheet.Columns.Item(1).ColumnWidth = 52.666666666666664
$newSheet.Columns.Item(2).ColumnWidth = 13.666666666666666

# --- Header row ---
$newSheet.Range("A1").Value = "Tasks"
$newSheet.Range("B1").Value = "Responsible"
$newSheet.Range("C1").Value = "Status"
$newSheet.Range("D1").Value = 1
$newSheet.Range("E1").Value = 2
$newSheet.Range("F1").Value = 3
$newSheet.Range("G1").Value = 4
$newSheet.Range("A1:G1").Font.Bold = $true
$newSheet.Range("A1:G1").WrapText = $true

# --- M7 task rows (written in this order to line up new shared-string ids) ---
$newSheet.Range("A2").Value = "Class Diagram"
$newSheet.Range("B2").Value = "Bhavesh"
$newSheet.Range("A3").Value = "Save/Load State (UI and controller and everything)"
$newSheet.Range("B3").Value = "Hunter"
$newSheet.Range("B4").Value = "Naman/ Pranil"
$newSheet.Range("A5").Value = "Code Critique and Java Doc"
$newSheet.Range("B5").Value = "Stephen"
$newSheet.Range("A4").Value = "Random event generation and implementation and extra credit"

# --- Page margins matching the other sprint sheets ---
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# --- Make the new sheet the active tab, zoomed like the rest, selection at A6 ---
$newSheet.Activate()
$excel.ActiveWindow.Zoom = 125
$newSheet.Range("A6").Select()
